$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4566
$ws.Range("I15").Value = 4566
$ws.Range("K15").Value = 13698
$ws.Range("M15").Value = -13529
$ws.Range("H18").Value = 1090
$ws.Range("I18").Value = 1250
$ws.Range("J18").Value = 450
$ws.Range("K18").Value = 1250
$ws.Range("L18").Value = 450
$ws.Range("M18").Value = -966
$ws.Range("N18").Value = -1018
$ws.Range("H116").Value = 3986.5
$ws.Range("I116").Value = 2172.25
$ws.Range("J116").Value = 4893.625
$ws.Range("K116").Value = 2172.25
$ws.Range("L116").Value = 4893.625
$ws.Range("M116").Value = 1269.75
$ws.Range("N116").Value = -11777.625
$ws.Range("H118").Value = 950
$ws.Range("I118").Value = 950
$ws.Range("K118").Value = 2850
$ws.Range("M118").Value = -1193
$ws.Range("H125").Value = 588128.75
$ws.Range("J125").Value = 5833
$ws.Range("L125").Value = 52497
$ws.Range("N125").Value = -57417
$ws.Range("H132").Value = 4643.8335
$ws.Range("I132").Value = 3342.2307
$ws.Range("J132").Value = 8028
$ws.Range("K132").Value = 10026.6921
$ws.Range("L132").Value = 24084
$ws.Range("M132").Value = -7496.6921
$ws.Range("N132").Value = -29144
$ws.Range("H135").Value = 440.54544
$ws.Range("I135").Value = 395.33334
$ws.Range("K135").Value = 3558.00006
$ws.Range("M135").Value = -1023.00006
$ws.Range("H137").Value = 2336.0588
$ws.Range("I137").Value = 2245.0435
$ws.Range("K137").Value = 6735.130500000001
$ws.Range("M137").Value = -4185.130500000001
$ws.Range("H138").Value = 2233.4546
$ws.Range("I138").Value = 1547.2307
$ws.Range("J138").Value = 3224.6667
$ws.Range("K138").Value = 4641.6921
$ws.Range("L138").Value = 9674.000100000001
$ws.Range("M138").Value = 498.3078999999998
$ws.Range("N138").Value = -19954.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2064.5693
$ws.Range("I32").Value = 2064.5693
$ws.Range("K32").Value = 2064.5693
$ws.Range("M32").Value = -1777.5693
$ws.Range("H57").Value = 9199.5
$ws.Range("I57").Value = 9199.5
$ws.Range("K57").Value = 9199.5
$ws.Range("M57").Value = -8715.5
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H96").Value = 92333.336
$ws.Range("J96").Value = 92333.336
$ws.Range("L96").Value = 92333.336
$ws.Range("N96").Value = -97825.336
$ws.Range("H132").Value = 2593.2896
$ws.Range("I132").Value = 2528.2432
$ws.Range("K132").Value = 7584.7296
$ws.Range("M132").Value = -5054.7296

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 3000
$ws.Range("I128").Value = 3000
$ws.Range("K128").Value = 9000
$ws.Range("M128").Value = -6510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7864.615
$ws.Range("I31").Value = 6999.6665
$ws.Range("K31").Value = 6999.6665
$ws.Range("M31").Value = -6704.6665
$ws.Range("H34").Value = 7864.615
$ws.Range("I34").Value = 6999.6665
$ws.Range("K34").Value = 6999.6665
$ws.Range("M34").Value = -6797.6665
$ws.Range("H138").Value = 333331.66
$ws.Range("I138").Value = 200000
$ws.Range("K138").Value = 200000
$ws.Range("M138").Value = -194860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H38").Value = 868.1539
$ws.Range("I38").Value = 492.75
$ws.Range("J38").Value = 1035
$ws.Range("K38").Value = 1478.25
$ws.Range("L38").Value = 3105
$ws.Range("M38").Value = -1131.25
$ws.Range("N38").Value = -3799
$ws.Range("H80").Value = 3999.75
$ws.Range("J80").Value = 3999.75
$ws.Range("L80").Value = 11999.25
$ws.Range("N80").Value = -13871.25
$ws.Range("H83").Value = 3999.75
$ws.Range("J83").Value = 3999.75
$ws.Range("L83").Value = 35997.75
$ws.Range("N83").Value = -45357.75
$ws.Range("H125").Value = 20000
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H131").Value = 19232928
$ws.Range("I131").Value = 125000950
$ws.Range("J131").Value = 2377.9546
$ws.Range("K131").Value = 375002850
$ws.Range("L131").Value = 7133.8638
$ws.Range("M131").Value = -374997810
$ws.Range("N131").Value = -17213.8638
$ws.Range("H132").Value = 2649.818
$ws.Range("I132").Value = 1981
$ws.Range("K132").Value = 17829
$ws.Range("M132").Value = -15299
$ws.Range("H139").Value = 3318.75
$ws.Range("I139").Value = 3259.3333
$ws.Range("J139").Value = 3497
$ws.Range("K139").Value = 9777.999899999999
$ws.Range("L139").Value = 10491
$ws.Range("M139").Value = -4637.999899999999
$ws.Range("N139").Value = -20771

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 137.28
$ws.Range("I2").Value = 104.588234
$ws.Range("K2").Value = 104.588234
$ws.Range("M2").Value = 8.411766
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1628.6
$ws.Range("I16").Value = 1753.8
$ws.Range("J16").Value = 1378.2
$ws.Range("K16").Value = 1753.8
$ws.Range("L16").Value = 1378.2
$ws.Range("M16").Value = -1583.8
$ws.Range("N16").Value = -1718.2
$ws.Range("H22").Value = 1564.4546
$ws.Range("I22").Value = 1518.4286
$ws.Range("J22").Value = 1645
$ws.Range("K22").Value = 1518.4286
$ws.Range("L22").Value = 1645
$ws.Range("M22").Value = -1223.4286
$ws.Range("N22").Value = -2235
$ws.Range("H27").Value = 1564.4546
$ws.Range("I27").Value = 1518.4286
$ws.Range("J27").Value = 1645
$ws.Range("K27").Value = 1518.4286
$ws.Range("L27").Value = 1645
$ws.Range("M27").Value = -1411.4286
$ws.Range("N27").Value = -1859
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H61").Value = 113440.445
$ws.Range("I61").Value = 127370.5
$ws.Range("K61").Value = 127370.5
$ws.Range("M61").Value = -127168.5
$ws.Range("H68").Value = 2328.2856
$ws.Range("J68").Value = 2499.3333
$ws.Range("L68").Value = 2499.3333
$ws.Range("N68").Value = -3997.3333
$ws.Range("H71").Value = 2328.2856
$ws.Range("J71").Value = 2499.3333
$ws.Range("L71").Value = 12496.6665
$ws.Range("N71").Value = -19984.6665
$ws.Range("H113").Value = 113440.445
$ws.Range("I113").Value = 127370.5
$ws.Range("K113").Value = 127370.5
$ws.Range("M113").Value = -125200.5
$ws.Range("H122").Value = 4521.091
$ws.Range("I122").Value = 4562.375
$ws.Range("J122").Value = 4411
$ws.Range("K122").Value = 13687.125
$ws.Range("L122").Value = 13233
$ws.Range("M122").Value = -11237.125
$ws.Range("N122").Value = -18133
$ws.Range("H132").Value = 8617.468000000001
$ws.Range("I132").Value = 8675
$ws.Range("J132").Value = 8377.75
$ws.Range("K132").Value = 26025
$ws.Range("L132").Value = 25133.25
$ws.Range("M132").Value = -23495
$ws.Range("N132").Value = -30193.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 103999.5
$ws.Range("J46").Value = 103999.5
$ws.Range("L46").Value = 103999.5
$ws.Range("N46").Value = -104461.5
$ws.Range("H104").Value = 30663.334
$ws.Range("I104").Value = 11990
$ws.Range("J104").Value = 40000
$ws.Range("K104").Value = 11990
$ws.Range("L104").Value = 40000
$ws.Range("M104").Value = -8496
$ws.Range("N104").Value = -46988
$ws.Range("H132").Value = 3152.4167
$ws.Range("I132").Value = 3152.4167
$ws.Range("K132").Value = 9457.250100000001
$ws.Range("M132").Value = -6927.250100000001
$ws.Range("H134").Value = 103999.5
$ws.Range("J134").Value = 103999.5
$ws.Range("L134").Value = 311998.5
$ws.Range("N134").Value = -317068.5
